$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()
$ws.Range("H107").Value = 1001.52
$ws.Range("I107").Value = 1093.6923
$ws.Range("J107").Value = 901.6667
$ws.Range("K107").Value = 1093.6923
$ws.Range("L107").Value = 901.6667
$ws.Range("M107").Value = 826.3077000000001
$ws.Range("N107").Value = -4741.6667

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 283.94446
$ws.Range("I4").Value = 260.16666
$ws.Range("J4").Value = 295.83334
$ws.Range("K4").Value = 260.16666
$ws.Range("L4").Value = 295.83334
$ws.Range("M4").Value = -144.16666
$ws.Range("N4").Value = -527.83334
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H6").Value = 3000
$ws.Range("J6").Value = 3000
$ws.Range("L6").Value = 3000
$ws.Range("N6").Value = -3346
$ws.Range("H9").Value = 20000
$ws.Range("J9").Value = 20000
$ws.Range("L9").Value = 20000
$ws.Range("N9").Value = -20340
$ws.Range("H20").Value = 20000
$ws.Range("J20").Value = 20000
$ws.Range("L20").Value = 20000
$ws.Range("N20").Value = -20540
$ws.Range("H23").Value = 51000
$ws.Range("J23").Value = 51000
$ws.Range("L23").Value = 51000
$ws.Range("N23").Value = -51518
$ws.Range("H37").Value = 8790.4
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 8790.4
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 8790.4
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -9336.4
$ws.Range("H44").Value = 15731.75
$ws.Range("J44").Value = 15731.75
$ws.Range("L44").Value = 15731.75
$ws.Range("N44").Value = -16707.75
$ws.Range("H55").Value = 23610.375
$ws.Range("J55").Value = 23610.375
$ws.Range("L55").Value = 23610.375
$ws.Range("N55").Value = -24240.375
$ws.Range("H63").Value = 5252.273
$ws.Range("I63").Value = 1925
$ws.Range("J63").Value = 6500
$ws.Range("K63").Value = 1925
$ws.Range("L63").Value = 6500
$ws.Range("M63").Value = -1239
$ws.Range("N63").Value = -7872
$ws.Range("H66").Value = 5252.273
$ws.Range("I66").Value = 1925
$ws.Range("J66").Value = 6500
$ws.Range("K66").Value = 9625
$ws.Range("L66").Value = 32500
$ws.Range("M66").Value = -6193
$ws.Range("N66").Value = -39364
$ws.Range("H80").Value = 28082.5
$ws.Range("J80").Value = 28082.5
$ws.Range("L80").Value = 28082.5
$ws.Range("N80").Value = -30078.5
$ws.Range("H83").Value = 28082.5
$ws.Range("J83").Value = 28082.5
$ws.Range("L83").Value = 84247.5
$ws.Range("N83").Value = -94231.5
$ws.Range("H88").Value = 3777.6
$ws.Range("J88").Value = 2544
$ws.Range("L88").Value = 2544
$ws.Range("N88").Value = -3356
$ws.Range("H91").Value = 3777.6
$ws.Range("J91").Value = 2544
$ws.Range("L91").Value = 2544
$ws.Range("N91").Value = -5352
$ws.Range("H102").Value = 3472.5
$ws.Range("I102").Value = 2950
$ws.Range("J102").Value = 3995
$ws.Range("K102").Value = 2950
$ws.Range("L102").Value = 3995
$ws.Range("M102").Value = -1328
$ws.Range("N102").Value = -7239
$ws.Range("H110").Value = 1191.9131
$ws.Range("I110").Value = 1047.125
$ws.Range("K110").Value = 1047.125
$ws.Range("M110").Value = 997.875

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H15").Value = 7892
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 7892
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 7892
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -8346
$ws.Range("H35").Value = 19516
$ws.Range("J35").Value = 24774
$ws.Range("L35").Value = 24774
$ws.Range("N35").Value = -25394
$ws.Range("H82").Value = 16532.875
$ws.Range("I82").Value = 11660
$ws.Range("J82").Value = 19456.6
$ws.Range("K82").Value = 11660
$ws.Range("L82").Value = 19456.6
$ws.Range("M82").Value = -11277
$ws.Range("N82").Value = -20222.6
$ws.Range("H85").Value = 16532.875
$ws.Range("I85").Value = 11660
$ws.Range("J85").Value = 19456.6
$ws.Range("K85").Value = 11660
$ws.Range("L85").Value = 19456.6
$ws.Range("M85").Value = -10334
$ws.Range("N85").Value = -22108.6
$ws.Range("H86").Value = 1922.6207
$ws.Range("I86").Value = 1809.6666
$ws.Range("J86").Value = 2107.4546
$ws.Range("K86").Value = 1809.6666
$ws.Range("L86").Value = 2107.4546
$ws.Range("M86").Value = -686.6666
$ws.Range("N86").Value = -4353.4546
$ws.Range("H89").Value = 1922.6207
$ws.Range("I89").Value = 1809.6666
$ws.Range("J89").Value = 2107.4546
$ws.Range("K89").Value = 9048.333000000001
$ws.Range("L89").Value = 10537.273
$ws.Range("M89").Value = -3432.333000000001
$ws.Range("N89").Value = -21769.273

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 55.857143
$ws.Range("I7").Value = 56.125
$ws.Range("J7").Value = 55.5
$ws.Range("K7").Value = 56.125
$ws.Range("L7").Value = 55.5
$ws.Range("M7").Value = 56.875
$ws.Range("N7").Value = -281.5
$ws.Range("H31").Value = 1104.6061
$ws.Range("I31").Value = 1036.9231
$ws.Range("J31").Value = 1356
$ws.Range("K31").Value = 1036.9231
$ws.Range("L31").Value = 1356
$ws.Range("M31").Value = -741.9231
$ws.Range("N31").Value = -1946
$ws.Range("H34").Value = 1104.6061
$ws.Range("I34").Value = 1036.9231
$ws.Range("J34").Value = 1356
$ws.Range("K34").Value = 1036.9231
$ws.Range("L34").Value = 1356
$ws.Range("M34").Value = -834.9231
$ws.Range("N34").Value = -1760
$ws.Range("H58").Value = 2349.7188
$ws.Range("I58").Value = 1728.579
$ws.Range("J58").Value = 3257.5386
$ws.Range("K58").Value = 1728.579
$ws.Range("L58").Value = 3257.5386
$ws.Range("M58").Value = -1525.579
$ws.Range("N58").Value = -3663.5386
$ws.Range("H136").Value = 2349.7188
$ws.Range("I136").Value = 1728.579
$ws.Range("J136").Value = 3257.5386
$ws.Range("K136").Value = 5185.737
$ws.Range("L136").Value = 9772.6158
$ws.Range("M136").Value = -2635.737
$ws.Range("N136").Value = -14872.6158

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1252008
$ws.Range("J131").Value = 1564150.6
$ws.Range("L131").Value = 4692451.800000001
$ws.Range("N131").Value = -4702531.800000001

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H97").Value = 26533.334
$ws.Range("J97").Value = 26533.334
$ws.Range("L97").Value = 26533.334
$ws.Range("N97").Value = -28515.334
$ws.Range("H132").Value = 5168.606
$ws.Range("I132").Value = 5734.8184
$ws.Range("J132").Value = 4036.182
$ws.Range("K132").Value = 17204.4552
$ws.Range("L132").Value = 12108.546
$ws.Range("M132").Value = -14674.4552
$ws.Range("N132").Value = -17168.546
$ws.Range("H136").Value = 1432.5333
$ws.Range("I136").Value = 576
$ws.Range("J136").Value = 7000
$ws.Range("K136").Value = 1728
$ws.Range("L136").Value = 21000
$ws.Range("M136").Value = 822
$ws.Range("N136").Value = -26100

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3418.5
$ws.Range("I132").Value = 3876.742
$ws.Range("J132").Value = 2127.0908
$ws.Range("K132").Value = 11630.226
$ws.Range("L132").Value = 6381.2724
$ws.Range("M132").Value = -9100.226000000001
$ws.Range("N132").Value = -11441.2724
$ws.Range("H136").Value = 10859.4
$ws.Range("I136").Value = 11273.053
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 33819.159
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -31269.159
$ws.Range("N136").Value = -14100
